$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.500.07'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '3.980.65'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.61'
$ws.Range("E5").Value = '  +4.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.91'
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("D7").Value = '3.973.25'
$ws.Range("E7").Value = '  -1.58%  '
$ws.Range("E8").Value = '  -5.57%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.744'
$ws.Range("E10").Value = '  -3.80%  '
$ws.Range("E11").Value = '  -4.83%  '
$ws.Range("E12").Value = '  +18.22%  '
$ws.Range("E13").Value = '  -2.72%  '
$ws.Range("E14").Value = '  -3.12%  '
$ws.Range("D15").Value = '4.613.03'
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").Value = '3.977.59'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.97'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.59'
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("E20").Value = '  -2.71%  '
$ws.Range("D21").Value = '71.361.75'
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '429.06'
$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("E23").Value = '  +1.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '97.52'
$ws.Range("E24").Value = '  -6.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.22'
$ws.Range("E25").Value = '  +5.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.57'
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.52'
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("E28").Value = '  -2.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.75'
$ws.Range("E29").Value = '  +15.71%  '
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.69'
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.68'
$ws.Range("E32").Value = '  +12.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.29'
$ws.Range("E33").Value = '  +20.79%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.131'
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '692.53'
$ws.Range("E35").Value = '  +1.43%  '
$ws.Range("B36").Value = 'Cosmos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '13.43'
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '65.57'
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.438'
$ws.Range("E38").Value = '  +2.32%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.151'
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0822'
$ws.Range("E40").Value = '  -4.72%  '
$ws.Range("E41").Value = '  -3.01%  '
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.26'
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("E46").Value = '  -5.71%  '
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.74'
$ws.Range("E48").Value = '  +6.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.37'
$ws.Range("E49").Value = '  -4.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.01'
$ws.Range("E50").Value = '  -1.73%  '
$ws.Range("D51").Value = '2.816.53'
$ws.Range("E51").Value = '  +8.12%  '
